# qc workflow is completed and is being tested
#
# A new blank row is inserted at row 3 (pushing the existing rows 3-6 down
# to 4-7). The new row inherits its cell formatting from the row above it
# (row 2), matching Excel's default "Insert" behaviour, and is left empty
# of values. Selection is then moved to E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting cells down and formatting the new
# row from the row above (xlShiftDown = -4121, xlFormatFromLeftOrAbove = 0).
$ws.Rows(3).Insert(-4121, 0)

# Match the workbook's recorded selection after the edit.
$ws.Range("E3").Select()
